$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in A42 (was 1, now 3)
$ws.Range("A42").Value = 3

# Recalculate dependent formulas (A43, A64, etc.)
$excel.Calculate()

# Update the sheet view: scroll position and active selection
$ws.Activate()
$ws.Range("A42").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
